$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$days = @("Monday","Tuesday","Wednesday","Thursday","Friday","Saturday","Sunday")
$meals = @("Breakfast","Lunch","Dinner")

# Append "Week 2" block (rows 23-43) mirroring the existing "Week 1" block (rows 2-22)
$r = 23
for ($d = 0; $d -lt 7; $d++) {
    for ($m = 0; $m -lt 3; $m++) {
        if ($m -eq 0) {
            if ($d -eq 0) {
                $ws.Cells.Item($r, 1).Value = "Week 2"
            }
            $ws.Cells.Item($r, 2).Value = $days[$d]
        }
        $ws.Cells.Item($r, 3).Value = $meals[$m]
        $r = $r + 1
    }
}

# Grow Table1 to cover the newly added rows
$lo = $ws.ListObjects.Item(1)
[void]$lo.Resize($ws.Range("A1:F43"))

# Match the author's final selection
[void]$ws.Range("C44").Select()
